$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.210134744644165
$ws.Range("B1").Value = 2.809362173080444
$ws.Range("C1").Value = 8.738932609558105
$ws.Range("D1").Value = 2.024343967437744
$ws.Range("E1").Value = 1.133161544799805
